$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three rows whose Sending cluster = "MuSCs" (old rows 8-10) entirely.
$ws.Rows.Item(8).Resize(3).Delete() | Out-Null

# Remove the three rows whose Sending cluster = "ECs" (old rows 2-4) entirely;
# remaining FAPs/MuSCs rows shift up to rows 2-7.
$ws.Rows.Item(2).Resize(3).Delete() | Out-Null

# Recalculated TPM-derived metrics for the surviving six rows.
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Fgf7"
$ws.Range("C2").Value = "Nrp1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 17.39906333333333
$ws.Range("H2").Value = 52.19719000000001
$ws.Range("I2").Value = 0.9351306508759385
$ws.Range("J2").Value = 0.9351306508759385
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 123.2806423333333
$ws.Range("N2").Value = 369.841927
$ws.Range("O2").Value = 0.6241574062367528
$ws.Range("P2").Value = 0.6241574062367526
$ws.Range("Q2").Value = 2144.967703731681
$ws.Range("R2").Value = 19304.70933358513
$ws.Range("S2").Value = 0.5836687215432121
$ws.Range("T2").Value = 0.583668721543212

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Fgf7"
$ws.Range("C3").Value = "Nrp1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 17.39906333333333
$ws.Range("H3").Value = 52.19719000000001
$ws.Range("I3").Value = 0.9351306508759385
$ws.Range("J3").Value = 0.9351306508759385
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 47.26005833333333
$ws.Range("N3").Value = 141.780175
$ws.Range("O3").Value = 0.2392728888301323
$ws.Range("P3").Value = 0.2392728888301322
$ws.Range("Q3").Value = 822.2807480786944
$ws.Range("R3").Value = 7400.52673270825
$ws.Range("S3").Value = 0.2237514122686876
$ws.Range("T3").Value = 0.2237514122686876

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Fgf7"
$ws.Range("C4").Value = "Nrp1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 17.39906333333333
$ws.Range("H4").Value = 52.19719000000001
$ws.Range("I4").Value = 0.9351306508759385
$ws.Range("J4").Value = 0.9351306508759385
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 26.97460733333333
$ws.Range("N4").Value = 80.923822
$ws.Range("O4").Value = 0.136569704933115
$ws.Range("P4").Value = 0.136569704933115
$ws.Range("Q4").Value = 469.3329013844645
$ws.Range("R4").Value = 4223.99611246018
$ws.Range("S4").Value = 0.1277105170640387
$ws.Range("T4").Value = 0.1277105170640387

$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Fgf7"
$ws.Range("C5").Value = "Nrp1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.206960666666667
$ws.Range("H5").Value = 3.620882
$ws.Range("I5").Value = 0.06486934912406146
$ws.Range("J5").Value = 0.06486934912406146
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 123.2806423333333
$ws.Range("N5").Value = 369.841927
$ws.Range("O5").Value = 0.6241574062367528
$ws.Range("P5").Value = 0.6241574062367526
$ws.Range("Q5").Value = 148.7948862577349
$ws.Range("R5").Value = 1339.153976319614
$ws.Range("S5").Value = 0.04048868469354057
$ws.Range("T5").Value = 0.04048868469354056

$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Fgf7"
$ws.Range("C6").Value = "Nrp1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.206960666666667
$ws.Range("H6").Value = 3.620882
$ws.Range("I6").Value = 0.06486934912406146
$ws.Range("J6").Value = 0.06486934912406146
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 47.26005833333333
$ws.Range("N6").Value = 141.780175
$ws.Range("O6").Value = 0.2392728888301323
$ws.Range("P6").Value = 0.2392728888301322
$ws.Range("Q6").Value = 57.04103151270554
$ws.Range("R6").Value = 513.3692836143499
$ws.Range("S6").Value = 0.0155214765614446
$ws.Range("T6").Value = 0.01552147656144459

$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Fgf7"
$ws.Range("C7").Value = "Nrp1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.206960666666667
$ws.Range("H7").Value = 3.620882
$ws.Range("I7").Value = 0.06486934912406146
$ws.Range("J7").Value = 0.06486934912406146
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 26.97460733333333
$ws.Range("N7").Value = 80.923822
$ws.Range("O7").Value = 0.136569704933115
$ws.Range("P7").Value = 0.136569704933115
$ws.Range("Q7").Value = 32.55729005011155
$ws.Range("R7").Value = 293.015610451004
$ws.Range("S7").Value = 0.008859187869076298
$ws.Range("T7").Value = 0.008859187869076296

